$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear everything first so the shared-string table rebuilds fresh,
# letting us control the order new unique strings are interned in.
$ws.Cells.Clear()

# --- Seed shared strings in the exact desired interning order ---
$ws.Cells.Item(2, 2).Value = 'HKL'
$ws.Cells.Item(3, 2).Value = 'Spiral5'
$ws.Cells.Item(4, 2).Value = 'Holden'
$ws.Cells.Item(5, 2).Value = 'Rizzie Spiral'
$ws.Cells.Item(6, 2).Value = 'RotRing OmegaMax-90'
$ws.Cells.Item(7, 2).Value = 'Equal Angle'
$ws.Cells.Item(8, 2).Value = 'Tilt Rotate'
$ws.Cells.Item(9, 2).Value = 'CLR'
$ws.Cells.Item(10, 2).Value = 'Rizzie Hex'
$ws.Cells.Item(11, 2).Value = 'Matthies Hex'
$ws.Cells.Item(12, 2).Value = 'Tilt Rotate_Partial'
$ws.Cells.Item(13, 2).Value = 'RotRing OmegaMax-60'
$ws.Cells.Item(14, 2).Value = 'Equal Angle_Partial'
$ws.Cells.Item(15, 2).Value = 'Rizzie Hex_Partial'
$ws.Cells.Item(16, 2).Value = 'ND Single'
$ws.Cells.Item(17, 2).Value = 'RD Single'
$ws.Cells.Item(18, 2).Value = 'TD Single'
$ws.Cells.Item(19, 2).Value = 'Morris Single'
$ws.Cells.Item(20, 2).Value = 'Ring Perpendicular to ND'
$ws.Cells.Item(21, 2).Value = 'Ring Perpendicular to RD'
$ws.Cells.Item(22, 2).Value = 'Ring Perpendicular to TD'
$ws.Cells.Item(23, 2).Value = 'OffsetFTD'
$ws.Cells.Item(24, 2).Value = 'OffsetATD'
$ws.Cells.Item(25, 2).Value = 'OffsetF45'
$ws.Cells.Item(26, 2).Value = 'OffsetA45'
$ws.Cells.Item(27, 2).Value = 'OffsetFRD'
$ws.Cells.Item(28, 2).Value = 'OffsetARD'
$ws.Cells.Item(29, 2).Value = 'Gaussian Quadrature'
$ws.Cells.Item(30, 2).Value = 'Michael-CCHex'
$ws.Cells.Item(31, 2).Value = 'Michael-SNHex'
$ws.Cells.Item(2, 3).Value = '[5, 1, 1]'
$ws.Cells.Item(2, 4).Value = '[4, 2, 2]'
$ws.Cells.Item(2, 5).Value = '[3, 1, 1]'
$ws.Cells.Item(2, 6).Value = '[3, 3, 1]'
$ws.Cells.Item(2, 7).Value = '[2, 2, 2]'
$ws.Cells.Item(2, 8).Value = '[1, 1, 1]'
$ws.Cells.Item(2, 9).Value = '[3, 3, 3]'
$ws.Cells.Item(2, 10).Value = '[2, 2, 0]'
$ws.Cells.Item(2, 11).Value = '[2, 0, 0]'
$ws.Cells.Item(2, 12).Value = '[4, 0, 0]'
$ws.Cells.Item(2, 13).Value = '[4, 2, 0]'
$ws.Cells.Item(2, 14).Value = '1Pair-A'
$ws.Cells.Item(2, 15).Value = '1Pair-B'
$ws.Cells.Item(2, 16).Value = '2Pairs-A'
$ws.Cells.Item(2, 17).Value = '2Pairs-B'
$ws.Cells.Item(2, 18).Value = '3Pairs-A'
$ws.Cells.Item(2, 19).Value = '3Pairs-B'
$ws.Cells.Item(2, 20).Value = '3Pairs-C'
$ws.Cells.Item(2, 21).Value = '4Pairs'
$ws.Cells.Item(2, 22).Value = '5A4F'
$ws.Cells.Item(2, 23).Value = 'MaxUnique'

# --- Fill remaining cell values ---
$ws.Cells.Item(1, 2).Value = 0
$ws.Cells.Item(1, 3).Value = 1
$ws.Cells.Item(1, 4).Value = 2
$ws.Cells.Item(1, 5).Value = 3
$ws.Cells.Item(1, 6).Value = 4
$ws.Cells.Item(1, 7).Value = 5
$ws.Cells.Item(1, 8).Value = 6
$ws.Cells.Item(1, 9).Value = 7
$ws.Cells.Item(1, 10).Value = 8
$ws.Cells.Item(1, 11).Value = 9
$ws.Cells.Item(1, 12).Value = 10
$ws.Cells.Item(1, 13).Value = 11
$ws.Cells.Item(1, 14).Value = 12
$ws.Cells.Item(1, 15).Value = 13
$ws.Cells.Item(1, 16).Value = 14
$ws.Cells.Item(1, 17).Value = 15
$ws.Cells.Item(1, 18).Value = 16
$ws.Cells.Item(1, 19).Value = 17
$ws.Cells.Item(1, 20).Value = 18
$ws.Cells.Item(1, 21).Value = 19
$ws.Cells.Item(1, 22).Value = 20
$ws.Cells.Item(1, 23).Value = 21

$ws.Cells.Item(2, 1).Value = 0

$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 3).Value = 0.9970346581216269
$ws.Cells.Item(3, 4).Value = 1.002601104782952
$ws.Cells.Item(3, 5).Value = 0.9993433324860577
$ws.Cells.Item(3, 6).Value = 1.001483698478692
$ws.Cells.Item(3, 7).Value = 1.000146606255431
$ws.Cells.Item(3, 8).Value = 1.000146606255431
$ws.Cells.Item(3, 9).Value = 1.000146606255431
$ws.Cells.Item(3, 10).Value = 0.9991190920980851
$ws.Cells.Item(3, 11).Value = 0.9958488752456628
$ws.Cells.Item(3, 12).Value = 0.9958488752456628
$ws.Cells.Item(3, 13).Value = 0.9985468099014432
$ws.Cells.Item(3, 14).Value = 1.000146606255431
$ws.Cells.Item(3, 15).Value = 0.9991190920980851
$ws.Cells.Item(3, 16).Value = 0.997483983671874
$ws.Cells.Item(3, 17).Value = 0.9992312122920715
$ws.Cells.Item(3, 18).Value = 0.9983715245330597
$ws.Cells.Item(3, 19).Value = 0.9981037666099352
$ws.Cells.Item(3, 20).Value = 0.9983715245330597
$ws.Cells.Item(3, 21).Value = 0.9986144765213092
$ws.Cells.Item(3, 22).Value = 0.9989209024681337
$ws.Cells.Item(3, 23).Value = 0.9992655221712439

$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 3).Value = 0.8913938261148956
$ws.Cells.Item(4, 4).Value = 1.064403574456597
$ws.Cells.Item(4, 5).Value = 0.9805561914844311
$ws.Cells.Item(4, 6).Value = 1.073846154679323
$ws.Cells.Item(4, 7).Value = 0.9909301033071435
$ws.Cells.Item(4, 8).Value = 0.9909301033071435
$ws.Cells.Item(4, 9).Value = 0.9909301033071435
$ws.Cells.Item(4, 10).Value = 1.008163220043487
$ws.Cells.Item(4, 11).Value = 0.8452496542632355
$ws.Cells.Item(4, 12).Value = 0.8452496542632355
$ws.Cells.Item(4, 13).Value = 0.9608016424108297
$ws.Cells.Item(4, 14).Value = 0.9909301033071435
$ws.Cells.Item(4, 15).Value = 1.008163220043487
$ws.Cells.Item(4, 16).Value = 0.9267064371533613
$ws.Cells.Item(4, 17).Value = 0.994359705763959
$ws.Cells.Item(4, 18).Value = 0.9481143258712886
$ws.Cells.Item(4, 19).Value = 0.9446563552637178
$ws.Cells.Item(4, 20).Value = 0.9481143258712886
$ws.Cells.Item(4, 21).Value = 0.9562247922745742
$ws.Cells.Item(4, 22).Value = 0.963165854481088
$ws.Cells.Item(4, 23).Value = 0.9769180458449929

$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(5, 3).Value = 0.6063303622464423
$ws.Cells.Item(5, 4).Value = 0.8890131037890022
$ws.Cells.Item(5, 5).Value = 0.6743465203917202
$ws.Cells.Item(5, 6).Value = 1.328028182804382
$ws.Cells.Item(5, 7).Value = 1.449636481447839
$ws.Cells.Item(5, 8).Value = 1.449636481447839
$ws.Cells.Item(5, 9).Value = 1.449636481447839
$ws.Cells.Item(5, 10).Value = 1.146176879462884
$ws.Cells.Item(5, 11).Value = 0.2053338392415246
$ws.Cells.Item(5, 12).Value = 0.2053338392415246
$ws.Cells.Item(5, 13).Value = 1.11490885669172
$ws.Cells.Item(5, 14).Value = 1.449636481447839
$ws.Cells.Item(5, 15).Value = 1.146176879462884
$ws.Cells.Item(5, 16).Value = 0.6757553593522042
$ws.Cells.Item(5, 17).Value = 0.9102616999273019
$ws.Cells.Item(5, 18).Value = 0.9337157333840823
$ws.Cells.Item(5, 19).Value = 0.6752857463653762
$ws.Cells.Item(5, 20).Value = 0.9337157333840823
$ws.Cells.Item(5, 21).Value = 0.8688734301359918
$ws.Cells.Item(5, 22).Value = 0.9850260403983612
$ws.Cells.Item(5, 23).Value = 0.9267217782594392

$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(6, 3).Value = 0.8517588439382799
$ws.Cells.Item(6, 4).Value = 1.063384020310425
$ws.Cells.Item(6, 5).Value = 0.9399590064455908
$ws.Cells.Item(6, 6).Value = 1.115477640568907
$ws.Cells.Item(6, 7).Value = 1.121773016118908
$ws.Cells.Item(6, 8).Value = 1.121773016118908
$ws.Cells.Item(6, 9).Value = 1.121773016118908
$ws.Cells.Item(6, 10).Value = 1.020133141128997
$ws.Cells.Item(6, 11).Value = 0.8051326476675793
$ws.Cells.Item(6, 12).Value = 0.8051326476675793
$ws.Cells.Item(6, 13).Value = 0.9289826895838499
$ws.Cells.Item(6, 14).Value = 1.121773016118908
$ws.Cells.Item(6, 15).Value = 1.020133141128997
$ws.Cells.Item(6, 16).Value = 0.9126328943982882
$ws.Cells.Item(6, 17).Value = 0.980046073787294
$ws.Cells.Item(6, 18).Value = 0.9823462683051615
$ws.Cells.Item(6, 19).Value = 0.9217415984140557
$ws.Cells.Item(6, 20).Value = 0.9823462683051615
$ws.Cells.Item(6, 21).Value = 0.9717494528402688
$ws.Cells.Item(6, 22).Value = 1.001754165495997
$ws.Cells.Item(6, 23).Value = 0.9808251257203171

$ws.Cells.Item(7, 1).Value = 5
$ws.Cells.Item(7, 3).Value = 0.8616585779099439
$ws.Cells.Item(7, 4).Value = 1.082349973393372
$ws.Cells.Item(7, 5).Value = 0.9699575335230537
$ws.Cells.Item(7, 6).Value = 1.087822276700288
$ws.Cells.Item(7, 7).Value = 1.00716922834863
$ws.Cells.Item(7, 8).Value = 1.00716922834863
$ws.Cells.Item(7, 9).Value = 1.00716922834863
$ws.Cells.Item(7, 10).Value = 0.9909400529488483
$ws.Cells.Item(7, 11).Value = 0.8009174915544769
$ws.Cells.Item(7, 12).Value = 0.8009174915544769
$ws.Cells.Item(7, 13).Value = 0.9521607082146945
$ws.Cells.Item(7, 14).Value = 1.00716922834863
$ws.Cells.Item(7, 15).Value = 0.9909400529488483
$ws.Cells.Item(7, 16).Value = 0.8959287722516627
$ws.Cells.Item(7, 17).Value = 0.980448793235951
$ws.Cells.Item(7, 18).Value = 0.9330089242839851
$ws.Cells.Item(7, 19).Value = 0.920605026008793
$ws.Cells.Item(7, 20).Value = 0.9330089242839851
$ws.Cells.Item(7, 21).Value = 0.9422460765937523
$ws.Cells.Item(7, 22).Value = 0.9552307069447277
$ws.Cells.Item(7, 23).Value = 0.9691219803241634

$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 3).Value = 0.5078783235314903
$ws.Cells.Item(8, 4).Value = 1.378964038726094
$ws.Cells.Item(8, 5).Value = 0.9129820740257959
$ws.Cells.Item(8, 6).Value = 1.344692464690705
$ws.Cells.Item(8, 7).Value = 0.6890835318725309
$ws.Cells.Item(8, 8).Value = 0.6890835318725309
$ws.Cells.Item(8, 9).Value = 0.6890835318725309
$ws.Cells.Item(8, 10).Value = 0.8044690155992474
$ws.Cells.Item(8, 11).Value = 0.4237496077148585
$ws.Cells.Item(8, 12).Value = 0.4237496077148585
$ws.Cells.Item(8, 13).Value = 0.7362563411081109
$ws.Cells.Item(8, 14).Value = 0.6890835318725309
$ws.Cells.Item(8, 15).Value = 0.8044690155992474
$ws.Cells.Item(8, 16).Value = 0.614109311657053
$ws.Cells.Item(8, 17).Value = 0.8587255448125217
$ws.Cells.Item(8, 18).Value = 0.6391007183955456
$ws.Cells.Item(8, 19).Value = 0.7137335657799673
$ws.Cells.Item(8, 20).Value = 0.6391007183955456
$ws.Cells.Item(8, 21).Value = 0.7075710573031082
$ws.Cells.Item(8, 22).Value = 0.7038735522169928
$ws.Cells.Item(8, 23).Value = 0.849759424658604

$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 3).Value = 0.9872895600799918
$ws.Cells.Item(9, 4).Value = 1.002793800868001
$ws.Cells.Item(9, 5).Value = 0.9914925118934795
$ws.Cells.Item(9, 6).Value = 1.002065120720584
$ws.Cells.Item(9, 7).Value = 1.029595654090568
$ws.Cells.Item(9, 8).Value = 1.029595654090568
$ws.Cells.Item(9, 9).Value = 1.029595654090568
$ws.Cells.Item(9, 10).Value = 0.9898540267451656
$ws.Cells.Item(9, 11).Value = 0.9683125544788916
$ws.Cells.Item(9, 12).Value = 0.9683125544788916
$ws.Cells.Item(9, 13).Value = 1.004991481674402
$ws.Cells.Item(9, 14).Value = 1.029595654090568
$ws.Cells.Item(9, 15).Value = 0.9898540267451656
$ws.Cells.Item(9, 16).Value = 0.9790832906120286
$ws.Cells.Item(9, 17).Value = 0.9906732693193225
$ws.Cells.Item(9, 18).Value = 0.9959207451048749
$ws.Cells.Item(9, 19).Value = 0.9832196977058456
$ws.Cells.Item(9, 20).Value = 0.9959207451048749
$ws.Cells.Item(9, 21).Value = 0.9948136868020261
$ws.Cells.Item(9, 22).Value = 1.001770080259734
$ws.Cells.Item(9, 23).Value = 0.9970493388188854

$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 3).Value = 0.9965388071352881
$ws.Cells.Item(10, 4).Value = 1.001487231456273
$ws.Cells.Item(10, 5).Value = 1.000946005549058
$ws.Cells.Item(10, 6).Value = 0.9976659926894591
$ws.Cells.Item(10, 7).Value = 1.000253679510076
$ws.Cells.Item(10, 8).Value = 1.000253679510076
$ws.Cells.Item(10, 9).Value = 1.000253679510076
$ws.Cells.Item(10, 10).Value = 1.001837648753831
$ws.Cells.Item(10, 11).Value = 0.9992657272407295
$ws.Cells.Item(10, 12).Value = 0.9992657272407295
$ws.Cells.Item(10, 13).Value = 1.003217031350375
$ws.Cells.Item(10, 14).Value = 1.000253679510076
$ws.Cells.Item(10, 15).Value = 1.001837648753831
$ws.Cells.Item(10, 16).Value = 1.00055168799728
$ws.Cells.Item(10, 17).Value = 1.001391827151445
$ws.Cells.Item(10, 18).Value = 1.000452351834879
$ws.Cells.Item(10, 19).Value = 1.000683127181206
$ws.Cells.Item(10, 20).Value = 1.000452351834879
$ws.Cells.Item(10, 21).Value = 1.000575765263424
$ws.Cells.Item(10, 22).Value = 1.000511348112754
$ws.Cells.Item(10, 23).Value = 1.000151515460636

$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 3).Value = 0.9762117737015955
$ws.Cells.Item(11, 4).Value = 1.004086679286821
$ws.Cells.Item(11, 5).Value = 0.9852569518642905
$ws.Cells.Item(11, 6).Value = 1.004304202202002
$ws.Cells.Item(11, 7).Value = 1.039481712045323
$ws.Cells.Item(11, 8).Value = 1.039481712045323
$ws.Cells.Item(11, 9).Value = 1.039481712045323
$ws.Cells.Item(11, 10).Value = 0.9835731654808404
$ws.Cells.Item(11, 11).Value = 0.9459411958895237
$ws.Cells.Item(11, 12).Value = 0.9459411958895237
$ws.Cells.Item(11, 13).Value = 1.012129349894452
$ws.Cells.Item(11, 14).Value = 1.039481712045323
$ws.Cells.Item(11, 15).Value = 0.9835731654808404
$ws.Cells.Item(11, 16).Value = 0.964757180685182
$ws.Cells.Item(11, 17).Value = 0.9844150586725655
$ws.Cells.Item(11, 18).Value = 0.9896653578052289
$ws.Cells.Item(11, 19).Value = 0.9715904377448848
$ws.Cells.Item(11, 20).Value = 0.9896653578052289
$ws.Cells.Item(11, 21).Value = 0.9885632563199943
$ws.Cells.Item(11, 22).Value = 0.9987469474650599
$ws.Cells.Item(11, 23).Value = 0.993873128795606

$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 3).Value = 0.4997076018098052
$ws.Cells.Item(12, 4).Value = 1.388271883793507
$ws.Cells.Item(12, 5).Value = 0.9183244832185233
$ws.Cells.Item(12, 6).Value = 1.346634258283066
$ws.Cells.Item(12, 7).Value = 0.6788853167362303
$ws.Cells.Item(12, 8).Value = 0.6788853167362303
$ws.Cells.Item(12, 9).Value = 0.6788853167362303
$ws.Cells.Item(12, 10).Value = 0.7947546501821317
$ws.Cells.Item(12, 11).Value = 0.4034084581922302
$ws.Cells.Item(12, 12).Value = 0.4034084581922302
$ws.Cells.Item(12, 13).Value = 0.7366351770368206
$ws.Cells.Item(12, 14).Value = 0.6788853167362303
$ws.Cells.Item(12, 15).Value = 0.7947546501821317
$ws.Cells.Item(12, 16).Value = 0.5990815541871809
$ws.Cells.Item(12, 17).Value = 0.8565395667003275
$ws.Cells.Item(12, 18).Value = 0.6256828083701974
$ws.Cells.Item(12, 19).Value = 0.7054958638642951
$ws.Cells.Item(12, 20).Value = 0.6256828083701974
$ws.Cells.Item(12, 21).Value = 0.6988432270822789
$ws.Cells.Item(12, 22).Value = 0.6948516450130692
$ws.Cells.Item(12, 23).Value = 0.8458277286565392

$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 3).Value = 0.9291251859951412
$ws.Cells.Item(13, 4).Value = 1.077245649343046
$ws.Cells.Item(13, 5).Value = 0.9889900242541831
$ws.Cells.Item(13, 6).Value = 1.103048546339144
$ws.Cells.Item(13, 7).Value = 0.9518236976823181
$ws.Cells.Item(13, 8).Value = 0.9518236976823181
$ws.Cells.Item(13, 9).Value = 0.9518236976823181
$ws.Cells.Item(13, 10).Value = 1.012088944983497
$ws.Cells.Item(13, 11).Value = 0.9570621675747396
$ws.Cells.Item(13, 12).Value = 0.9570621675747396
$ws.Cells.Item(13, 13).Value = 0.9008141409332754
$ws.Cells.Item(13, 14).Value = 0.9518236976823181
$ws.Cells.Item(13, 15).Value = 1.012088944983497
$ws.Cells.Item(13, 16).Value = 0.9845755562791181
$ws.Cells.Item(13, 17).Value = 1.00053948461884
$ws.Cells.Item(13, 18).Value = 0.9736582700801848
$ws.Cells.Item(13, 19).Value = 0.9860470456041398
$ws.Cells.Item(13, 20).Value = 0.9736582700801848
$ws.Cells.Item(13, 21).Value = 0.9774912086236844
$ws.Cells.Item(13, 22).Value = 0.9723577064354112
$ws.Cells.Item(13, 23).Value = 0.9900247946381681

$ws.Cells.Item(14, 1).Value = 12
$ws.Cells.Item(14, 3).Value = 0.8531764907578956
$ws.Cells.Item(14, 4).Value = 1.092026461073684
$ws.Cells.Item(14, 5).Value = 0.9801807238105261
$ws.Cells.Item(14, 6).Value = 1.088435167568421
$ws.Cells.Item(14, 7).Value = 0.9636229743053657
$ws.Cells.Item(14, 8).Value = 0.9636229743053657
$ws.Cells.Item(14, 9).Value = 0.9636229743053657
$ws.Cells.Item(14, 10).Value = 0.9688389514473683
$ws.Cells.Item(14, 11).Value = 0.7838774836560001
$ws.Cells.Item(14, 12).Value = 0.7838774836560001
$ws.Cells.Item(14, 13).Value = 0.9582893362421041
$ws.Cells.Item(14, 14).Value = 0.9636229743053657
$ws.Cells.Item(14, 15).Value = 0.9688389514473683
$ws.Cells.Item(14, 16).Value = 0.8763582175516842
$ws.Cells.Item(14, 17).Value = 0.9745098376289472
$ws.Cells.Item(14, 18).Value = 0.9054464698029113
$ws.Cells.Item(14, 19).Value = 0.9109657196379648
$ws.Cells.Item(14, 20).Value = 0.9054464698029113
$ws.Cells.Item(14, 21).Value = 0.924130033304815
$ws.Cells.Item(14, 22).Value = 0.932028621504925
$ws.Cells.Item(14, 23).Value = 0.9610559486076706

$ws.Cells.Item(15, 1).Value = 13
$ws.Cells.Item(15, 3).Value = 1.084802360206074
$ws.Cells.Item(15, 4).Value = 0.9679369641928526
$ws.Cells.Item(15, 5).Value = 1.045716550664537
$ws.Cells.Item(15, 6).Value = 0.9378020797330487
$ws.Cells.Item(15, 7).Value = 0.9444415930157244
$ws.Cells.Item(15, 8).Value = 0.9444415930157244
$ws.Cells.Item(15, 9).Value = 0.9444415930157244
$ws.Cells.Item(15, 10).Value = 0.9453598009131896
$ws.Cells.Item(15, 11).Value = 1.07154081820489
$ws.Cells.Item(15, 12).Value = 1.07154081820489
$ws.Cells.Item(15, 13).Value = 1.04353551595399
$ws.Cells.Item(15, 14).Value = 0.9444415930157244
$ws.Cells.Item(15, 15).Value = 0.9453598009131896
$ws.Cells.Item(15, 16).Value = 1.00845030955904
$ws.Cells.Item(15, 17).Value = 0.9955381757888635
$ws.Cells.Item(15, 18).Value = 0.9871140707112679
$ws.Cells.Item(15, 19).Value = 1.020872389927539
$ws.Cells.Item(15, 20).Value = 0.987114070711268
$ws.Cells.Item(15, 21).Value = 1.001764690699585
$ws.Cells.Item(15, 22).Value = 0.9903000711628132
$ws.Cells.Item(15, 23).Value = 1.005141960360538

$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 3).Value = 0.1250664500000001
$ws.Cells.Item(16, 4).Value = 1.690899200000001
$ws.Cells.Item(16, 5).Value = 0.8444181999999999
$ws.Cells.Item(16, 6).Value = 1.618994999999999
$ws.Cells.Item(16, 7).Value = 0.3721459699999993
$ws.Cells.Item(16, 8).Value = 0.3721459699999993
$ws.Cells.Item(16, 9).Value = 0.3721459699999993
$ws.Cells.Item(16, 10).Value = 0.60377207
$ws.Cells.Item(16, 11).Value = 0.002395270000000003
$ws.Cells.Item(16, 12).Value = 0.002395270000000003
$ws.Cells.Item(16, 13).Value = 0.5142619599999995
$ws.Cells.Item(16, 14).Value = 0.3721459699999993
$ws.Cells.Item(16, 15).Value = 0.60377207
$ws.Cells.Item(16, 16).Value = 0.30308367
$ws.Cells.Item(16, 17).Value = 0.724095135
$ws.Cells.Item(16, 18).Value = 0.3261044366666664
$ws.Cells.Item(16, 19).Value = 0.4835285133333333
$ws.Cells.Item(16, 20).Value = 0.3261044366666664
$ws.Cells.Item(16, 21).Value = 0.4556828774999998
$ws.Cells.Item(16, 22).Value = 0.4389754959999997
$ws.Cells.Item(16, 23).Value = 0.7214942649999998

$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 3).Value = 0.062230451
$ws.Cells.Item(17, 4).Value = 2.826607699999999
$ws.Cells.Item(17, 5).Value = 0.7729728300000001
$ws.Cells.Item(17, 6).Value = 0.99995383
$ws.Cells.Item(17, 7).Value = 1.913908
$ws.Cells.Item(17, 8).Value = 1.913908
$ws.Cells.Item(17, 9).Value = 1.913908
$ws.Cells.Item(17, 10).Value = 0.09556448200000001
$ws.Cells.Item(17, 11).Value = -0.0012267103
$ws.Cells.Item(17, 12).Value = -0.0012267103
$ws.Cells.Item(17, 13).Value = 0.08389542799999999
$ws.Cells.Item(17, 14).Value = 1.913908
$ws.Cells.Item(17, 15).Value = 0.09556448200000001
$ws.Cells.Item(17, 16).Value = 0.04716888585
$ws.Cells.Item(17, 17).Value = 0.434268656
$ws.Cells.Item(17, 18).Value = 0.6694152572333333
$ws.Cells.Item(17, 19).Value = 0.2891035339
$ws.Cells.Item(17, 20).Value = 0.6694152572333333
$ws.Cells.Item(17, 21).Value = 0.695304650425
$ws.Cells.Item(17, 22).Value = 0.9390253203400001
$ws.Cells.Item(17, 23).Value = 0.8442382513375

$ws.Cells.Item(18, 1).Value = 16
$ws.Cells.Item(18, 3).Value = 0.0289053
$ws.Cells.Item(18, 4).Value = 1.4482004
$ws.Cells.Item(18, 5).Value = 0.37231391
$ws.Cells.Item(18, 6).Value = 2.3464701
$ws.Cells.Item(18, 7).Value = 1.0627484
$ws.Cells.Item(18, 8).Value = 1.0627484
$ws.Cells.Item(18, 9).Value = 1.0627484
$ws.Cells.Item(18, 10).Value = 0.57031013
$ws.Cells.Item(18, 11).Value = -0.0026079641
$ws.Cells.Item(18, 12).Value = -0.0026079641
$ws.Cells.Item(18, 13).Value = 0.18297177
$ws.Cells.Item(18, 14).Value = 1.0627484
$ws.Cells.Item(18, 15).Value = 0.57031013
$ws.Cells.Item(18, 16).Value = 0.28385108295
$ws.Cells.Item(18, 17).Value = 0.47131202
$ws.Cells.Item(18, 18).Value = 0.5434835219666666
$ws.Cells.Item(18, 19).Value = 0.3133386919666667
$ws.Cells.Item(18, 20).Value = 0.5434835219666666
$ws.Cells.Item(18, 21).Value = 0.500691118975
$ws.Cells.Item(18, 22).Value = 0.6131025751799999
$ws.Cells.Item(18, 23).Value = 0.7511640057375

$ws.Cells.Item(19, 1).Value = 17
$ws.Cells.Item(19, 3).Value = 0.90363594
$ws.Cells.Item(19, 4).Value = 0.7938472600000001
$ws.Cells.Item(19, 5).Value = 0.9695692299999998
$ws.Cells.Item(19, 6).Value = 1.0837992
$ws.Cells.Item(19, 7).Value = 1.881306
$ws.Cells.Item(19, 8).Value = 1.881306
$ws.Cells.Item(19, 9).Value = 1.881306
$ws.Cells.Item(19, 10).Value = 0.16949639
$ws.Cells.Item(19, 11).Value = 0.18124095
$ws.Cells.Item(19, 12).Value = 0.18124095
$ws.Cells.Item(19, 13).Value = 1.1129315
$ws.Cells.Item(19, 14).Value = 1.881306
$ws.Cells.Item(19, 15).Value = 0.16949639
$ws.Cells.Item(19, 16).Value = 0.17536867
$ws.Cells.Item(19, 17).Value = 0.5695328099999999
$ws.Cells.Item(19, 18).Value = 0.7440144466666666
$ws.Cells.Item(19, 19).Value = 0.4401021899999999
$ws.Cells.Item(19, 20).Value = 0.7440144466666666
$ws.Cells.Item(19, 21).Value = 0.8004031424999999
$ws.Cells.Item(19, 22).Value = 1.016583714
$ws.Cells.Item(19, 23).Value = 0.8869783087499999

$ws.Cells.Item(20, 1).Value = 18
$ws.Cells.Item(20, 3).Value = 0.7196994528904109
$ws.Cells.Item(20, 4).Value = 1.082589172328768
$ws.Cells.Item(20, 5).Value = 0.8295555247945203
$ws.Cells.Item(20, 6).Value = 1.035949045616438
$ws.Cells.Item(20, 7).Value = 1.508549973632877
$ws.Cells.Item(20, 8).Value = 1.508549973632877
$ws.Cells.Item(20, 9).Value = 1.508549973632877
$ws.Cells.Item(20, 10).Value = 0.8187743380547948
$ws.Cells.Item(20, 11).Value = 0.3408706750676713
$ws.Cells.Item(20, 12).Value = 0.3408706750676713
$ws.Cells.Item(20, 13).Value = 1.117080897315069
$ws.Cells.Item(20, 14).Value = 1.508549973632877
$ws.Cells.Item(20, 15).Value = 0.8187743380547948
$ws.Cells.Item(20, 16).Value = 0.579822506561233
$ws.Cells.Item(20, 17).Value = 0.8241649314246575
$ws.Cells.Item(20, 18).Value = 0.8893983289184476
$ws.Cells.Item(20, 19).Value = 0.6630668459723288
$ws.Cells.Item(20, 20).Value = 0.8893983289184476
$ws.Cells.Item(20, 21).Value = 0.8744376278874658
$ws.Cells.Item(20, 22).Value = 1.001260097036548
$ws.Cells.Item(20, 23).Value = 0.9316336349625686

$ws.Cells.Item(21, 1).Value = 19
$ws.Cells.Item(21, 3).Value = 0.436316494263158
$ws.Cells.Item(21, 4).Value = 0.981641990526316
$ws.Cells.Item(21, 5).Value = 0.8766805452631579
$ws.Cells.Item(21, 6).Value = 1.271781517368421
$ws.Cells.Item(21, 7).Value = 1.225705933968421
$ws.Cells.Item(21, 8).Value = 1.225705933968421
$ws.Cells.Item(21, 9).Value = 1.225705933968421
$ws.Cells.Item(21, 10).Value = 1.302075806315789
$ws.Cells.Item(21, 11).Value = 0.05247285255052631
$ws.Cells.Item(21, 12).Value = 0.05247285255052631
$ws.Cells.Item(21, 13).Value = 1.150126871578947
$ws.Cells.Item(21, 14).Value = 1.225705933968421
$ws.Cells.Item(21, 15).Value = 1.302075806315789
$ws.Cells.Item(21, 16).Value = 0.6772743294331578
$ws.Cells.Item(21, 17).Value = 1.089378175789474
$ws.Cells.Item(21, 18).Value = 0.8600848642782456
$ws.Cells.Item(21, 19).Value = 0.7437430680431577
$ws.Cells.Item(21, 20).Value = 0.8600848642782456
$ws.Cells.Item(21, 21).Value = 0.8642337845244736
$ws.Cells.Item(21, 22).Value = 0.9365282144132632
$ws.Cells.Item(21, 23).Value = 0.9121002514793422

$ws.Cells.Item(22, 1).Value = 20
$ws.Cells.Item(22, 3).Value = 0.5610245528947368
$ws.Cells.Item(22, 4).Value = 0.9849497036842106
$ws.Cells.Item(22, 5).Value = 0.6876247257894736
$ws.Cells.Item(22, 6).Value = 1.336754370526316
$ws.Cells.Item(22, 7).Value = 1.393779560489474
$ws.Cells.Item(22, 8).Value = 1.393779560489474
$ws.Cells.Item(22, 9).Value = 1.393779560489474
$ws.Cells.Item(22, 10).Value = 1.080590806368421
$ws.Cells.Item(22, 11).Value = 0.1862143961552631
$ws.Cells.Item(22, 12).Value = 0.1862143961552631
$ws.Cells.Item(22, 13).Value = 1.054505948105263
$ws.Cells.Item(22, 14).Value = 1.393779560489474
$ws.Cells.Item(22, 15).Value = 1.080590806368421
$ws.Cells.Item(22, 16).Value = 0.6334026012618421
$ws.Cells.Item(22, 17).Value = 0.8841077660789474
$ws.Cells.Item(22, 18).Value = 0.8868615876710527
$ws.Cells.Item(22, 19).Value = 0.6514766427710527
$ws.Cells.Item(22, 20).Value = 0.8868615876710527
$ws.Cells.Item(22, 21).Value = 0.837052372200658
$ws.Cells.Item(22, 22).Value = 0.9483978098584211
$ws.Cells.Item(22, 23).Value = 0.9106805080016449

$ws.Cells.Item(23, 1).Value = 21
$ws.Cells.Item(23, 3).Value = 0.7630427675592577
$ws.Cells.Item(23, 4).Value = 0.683104443778473
$ws.Cells.Item(23, 5).Value = 0.9464104146673036
$ws.Cells.Item(23, 6).Value = 1.289716805531327
$ws.Cells.Item(23, 7).Value = 0.09994839146864848
$ws.Cells.Item(23, 8).Value = 0.09994839146864848
$ws.Cells.Item(23, 9).Value = 0.09994839146864848
$ws.Cells.Item(23, 10).Value = 1.976323055206932
$ws.Cells.Item(23, 11).Value = 0.244327500059894
$ws.Cells.Item(23, 12).Value = 0.244327500059894
$ws.Cells.Item(23, 13).Value = 1.407947286107057
$ws.Cells.Item(23, 14).Value = 0.09994839146864848
$ws.Cells.Item(23, 15).Value = 1.976323055206932
$ws.Cells.Item(23, 16).Value = 1.110325277633413
$ws.Cells.Item(23, 17).Value = 1.461366734937118
$ws.Cells.Item(23, 18).Value = 0.7735329822451581
$ws.Cells.Item(23, 19).Value = 1.055686989978043
$ws.Cells.Item(23, 20).Value = 0.7735329822451581
$ws.Cells.Item(23, 21).Value = 0.8167523403506944
$ws.Cells.Item(23, 22).Value = 0.6733915505742852
$ws.Cells.Item(23, 23).Value = 0.9263525830473616

$ws.Cells.Item(24, 1).Value = 22
$ws.Cells.Item(24, 3).Value = 1.126208870142665
$ws.Cells.Item(24, 4).Value = 1.199641628230355
$ws.Cells.Item(24, 5).Value = 1.059127154540698
$ws.Cells.Item(24, 6).Value = 0.9653365473281356
$ws.Cells.Item(24, 7).Value = 1.154819413131847
$ws.Cells.Item(24, 8).Value = 1.154819413131847
$ws.Cells.Item(24, 9).Value = 1.154819413131847
$ws.Cells.Item(24, 10).Value = 0.6158226443872522
$ws.Cells.Item(24, 11).Value = 1.487078587347189
$ws.Cells.Item(24, 12).Value = 1.487078587347189
$ws.Cells.Item(24, 13).Value = 0.6774639121706539
$ws.Cells.Item(24, 14).Value = 1.154819413131847
$ws.Cells.Item(24, 15).Value = 0.6158226443872522
$ws.Cells.Item(24, 16).Value = 1.051450615867221
$ws.Cells.Item(24, 17).Value = 0.8374748994639749
$ws.Cells.Item(24, 18).Value = 1.085906881622096
$ws.Cells.Item(24, 19).Value = 1.054009462091713
$ws.Cells.Item(24, 20).Value = 1.085906881622096
$ws.Cells.Item(24, 21).Value = 1.079211949851747
$ws.Cells.Item(24, 22).Value = 1.094333442507767
$ws.Cells.Item(24, 23).Value = 1.035687344659849

$ws.Cells.Item(25, 1).Value = 23
$ws.Cells.Item(25, 3).Value = 1.242818300640022
$ws.Cells.Item(25, 4).Value = 1.086102743082858
$ws.Cells.Item(25, 5).Value = 1.035393312944749
$ws.Cells.Item(25, 6).Value = 0.8421911193373977
$ws.Cells.Item(25, 7).Value = 1.11699080466401
$ws.Cells.Item(25, 8).Value = 1.11699080466401
$ws.Cells.Item(25, 9).Value = 1.11699080466401
$ws.Cells.Item(25, 10).Value = 0.6904590665431408
$ws.Cells.Item(25, 11).Value = 1.526267584329738
$ws.Cells.Item(25, 12).Value = 1.526267584329738
$ws.Cells.Item(25, 13).Value = 0.8507554068194179
$ws.Cells.Item(25, 14).Value = 1.11699080466401
$ws.Cells.Item(25, 15).Value = 0.6904590665431408
$ws.Cells.Item(25, 16).Value = 1.108363325436439
$ws.Cells.Item(25, 17).Value = 0.8629261897439446
$ws.Cells.Item(25, 18).Value = 1.111239151845629
$ws.Cells.Item(25, 19).Value = 1.084039987939209
$ws.Cells.Item(25, 20).Value = 1.111239151845629
$ws.Cells.Item(25, 21).Value = 1.092277692120409
$ws.Cells.Item(25, 22).Value = 1.09722031462913
$ws.Cells.Item(25, 23).Value = 1.048872292295167

$ws.Cells.Item(26, 1).Value = 24
$ws.Cells.Item(26, 3).Value = 1.010484817529409
$ws.Cells.Item(26, 4).Value = 0.8480083243235578
$ws.Cells.Item(26, 5).Value = 1.010380707912029
$ws.Cells.Item(26, 6).Value = 0.9351409646023666
$ws.Cells.Item(26, 7).Value = 1.003569769547607
$ws.Cells.Item(26, 8).Value = 1.003569769547607
$ws.Cells.Item(26, 9).Value = 1.003569769547607
$ws.Cells.Item(26, 10).Value = 1.053270246613401
$ws.Cells.Item(26, 11).Value = 0.8074113510347647
$ws.Cells.Item(26, 12).Value = 0.8074113510347647
$ws.Cells.Item(26, 13).Value = 1.217295016323791
$ws.Cells.Item(26, 14).Value = 1.003569769547607
$ws.Cells.Item(26, 15).Value = 1.053270246613401
$ws.Cells.Item(26, 16).Value = 0.9303407988240829
$ws.Cells.Item(26, 17).Value = 1.031825477262715
$ws.Cells.Item(26, 18).Value = 0.9547504557319243
$ws.Cells.Item(26, 19).Value = 0.9570207685200649
$ws.Cells.Item(26, 20).Value = 0.9547504557319243
$ws.Cells.Item(26, 21).Value = 0.9686580187769505
$ws.Cells.Item(26, 22).Value = 0.9756403689310819
$ws.Cells.Item(26, 23).Value = 0.9856951497358659

$ws.Cells.Item(27, 1).Value = 25
$ws.Cells.Item(27, 3).Value = 1.041136374355129
$ws.Cells.Item(27, 4).Value = 1.210556466762628
$ws.Cells.Item(27, 5).Value = 0.9768003691381658
$ws.Cells.Item(27, 6).Value = 0.786329835412645
$ws.Cells.Item(27, 7).Value = 1.578475872977647
$ws.Cells.Item(27, 8).Value = 1.578475872977647
$ws.Cells.Item(27, 9).Value = 1.578475872977647
$ws.Cells.Item(27, 10).Value = 0.5341849361534837
$ws.Cells.Item(27, 11).Value = 1.440449686404931
$ws.Cells.Item(27, 12).Value = 1.440449686404931
$ws.Cells.Item(27, 13).Value = 0.8389755295940976
$ws.Cells.Item(27, 14).Value = 1.578475872977647
$ws.Cells.Item(27, 15).Value = 0.5341849361534837
$ws.Cells.Item(27, 16).Value = 0.9873173112792071
$ws.Cells.Item(27, 17).Value = 0.7554926526458248
$ws.Cells.Item(27, 18).Value = 1.184370165178687
$ws.Cells.Item(27, 19).Value = 0.98381166389886
$ws.Cells.Item(27, 20).Value = 1.184370165178687
$ws.Cells.Item(27, 21).Value = 1.132477716168557
$ws.Cells.Item(27, 22).Value = 1.221677347530375
$ws.Cells.Item(27, 23).Value = 1.050863633849841

$ws.Cells.Item(28, 1).Value = 26
$ws.Cells.Item(28, 3).Value = 1.015387802006874
$ws.Cells.Item(28, 4).Value = 0.8880454304002154
$ws.Cells.Item(28, 5).Value = 1.046710292947735
$ws.Cells.Item(28, 6).Value = 1.059605845691872
$ws.Cells.Item(28, 7).Value = 0.7102682169571874
$ws.Cells.Item(28, 8).Value = 0.7102682169571874
$ws.Cells.Item(28, 9).Value = 0.7102682169571874
$ws.Cells.Item(28, 10).Value = 1.150125617032131
$ws.Cells.Item(28, 11).Value = 0.8745762512811038
$ws.Cells.Item(28, 12).Value = 0.8745762512811038
$ws.Cells.Item(28, 13).Value = 1.081311529761104
$ws.Cells.Item(28, 14).Value = 0.7102682169571874
$ws.Cells.Item(28, 15).Value = 1.150125617032131
$ws.Cells.Item(28, 16).Value = 1.012350934156617
$ws.Cells.Item(28, 17).Value = 1.098417954989933
$ws.Cells.Item(28, 18).Value = 0.9116566950901408
$ws.Cells.Item(28, 19).Value = 1.023804053753657
$ws.Cells.Item(28, 20).Value = 0.9116566950901408
$ws.Cells.Item(28, 21).Value = 0.9454200945545392
$ws.Cells.Item(28, 22).Value = 0.8983897190350689
$ws.Cells.Item(28, 23).Value = 0.9782538732597778

$ws.Cells.Item(29, 1).Value = 27
$ws.Cells.Item(29, 3).Value = 0.9601991765706627
$ws.Cells.Item(29, 4).Value = 0.923391638927975
$ws.Cells.Item(29, 5).Value = 0.8810148527434475
$ws.Cells.Item(29, 6).Value = 1.0300893289662
$ws.Cells.Item(29, 7).Value = 1.16721969902503
$ws.Cells.Item(29, 8).Value = 1.16721969902503
$ws.Cells.Item(29, 9).Value = 1.16721969902503
$ws.Cells.Item(29, 10).Value = 0.9248947266698714
$ws.Cells.Item(29, 11).Value = 1.139533613708962
$ws.Cells.Item(29, 12).Value = 1.139533613708962
$ws.Cells.Item(29, 13).Value = 1.028329068319012
$ws.Cells.Item(29, 14).Value = 1.16721969902503
$ws.Cells.Item(29, 15).Value = 0.9248947266698714
$ws.Cells.Item(29, 16).Value = 1.032214170189417
$ws.Cells.Item(29, 17).Value = 0.9029547897066594
$ws.Cells.Item(29, 18).Value = 1.077216013134621
$ws.Cells.Item(29, 19).Value = 0.9818143977074271
$ws.Cells.Item(29, 20).Value = 1.077216013134621
$ws.Cells.Item(29, 21).Value = 1.028165723036828
$ws.Cells.Item(29, 22).Value = 1.055976518234468
$ws.Cells.Item(29, 23).Value = 1.006834013116395

$ws.Cells.Item(30, 1).Value = 28
$ws.Cells.Item(30, 3).Value = 1.113589766783272
$ws.Cells.Item(30, 4).Value = 0.7952330968881779
$ws.Cells.Item(30, 5).Value = 1.033441328903443
$ws.Cells.Item(30, 6).Value = 1.212271154964766
$ws.Cells.Item(30, 7).Value = 1.196543162745848
$ws.Cells.Item(30, 8).Value = 1.196543162745848
$ws.Cells.Item(30, 9).Value = 1.196543162745848
$ws.Cells.Item(30, 10).Value = 0.9136697109278663
$ws.Cells.Item(30, 11).Value = 0.7707884210560237
$ws.Cells.Item(30, 12).Value = 0.7707884210560237
$ws.Cells.Item(30, 13).Value = 0.9317898461780962
$ws.Cells.Item(30, 14).Value = 1.196543162745848
$ws.Cells.Item(30, 15).Value = 0.9136697109278663
$ws.Cells.Item(30, 16).Value = 0.842229065991945
$ws.Cells.Item(30, 17).Value = 0.9735555199156545
$ws.Cells.Item(30, 18).Value = 0.9603337649099126
$ws.Cells.Item(30, 19).Value = 0.9059664869624443
$ws.Cells.Item(30, 20).Value = 0.9603337649099126
$ws.Cells.Item(30, 21).Value = 0.9786106559082952
$ws.Cells.Item(30, 22).Value = 1.022197157275806
$ws.Cells.Item(30, 23).Value = 0.9959158110559365

$ws.Cells.Item(31, 1).Value = 29
$ws.Cells.Item(31, 3).Value = 1.434933761042172
$ws.Cells.Item(31, 4).Value = 0.7185817355205814
$ws.Cells.Item(31, 5).Value = 1.272992228463027
$ws.Cells.Item(31, 6).Value = 1.223818482429671
$ws.Cells.Item(31, 7).Value = 1.085414146119072
$ws.Cells.Item(31, 8).Value = 1.085414146119072
$ws.Cells.Item(31, 9).Value = 1.085414146119072
$ws.Cells.Item(31, 10).Value = 0.5368864095795753
$ws.Cells.Item(31, 11).Value = 0.5613535121170528
$ws.Cells.Item(31, 12).Value = 0.5613535121170528
$ws.Cells.Item(31, 13).Value = 0.8995686629137188
$ws.Cells.Item(31, 14).Value = 1.085414146119072
$ws.Cells.Item(31, 15).Value = 0.5368864095795753
$ws.Cells.Item(31, 16).Value = 0.549119960848314
$ws.Cells.Item(31, 17).Value = 0.9049393190213011
$ws.Cells.Item(31, 18).Value = 0.7278846892719001
$ws.Cells.Item(31, 19).Value = 0.7904107167198849
$ws.Cells.Item(31, 20).Value = 0.7278846892719001
$ws.Cells.Item(31, 21).Value = 0.8641615740696817
$ws.Cells.Item(31, 22).Value = 0.9084120884795599
$ws.Cells.Item(31, 23).Value = 0.9666936172731089

# --- Copy formatting for newly added A30/A31 cells ---
$ws.Range("A29").Copy() | Out-Null
$ws.Range("A30:A31").PasteSpecial(-4122)
$excel.CutCopyMode = $false
